$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the *2023 column (H) data values for rows 2-11 ---
$ws.Range("H2").Value = 21.6
$ws.Range("H3").Value = 21.6
$ws.Range("H4").Value = 21.6
$ws.Range("H5").Value = 21.6
$ws.Range("H6").Value = 21.6
$ws.Range("H7").Value = 21.6
$ws.Range("H8").Value = 21.6
$ws.Range("H9").Value = 21.6
$ws.Range("H10").Value = 21.6
$ws.Range("H11").Value = 21.6

# --- Totals row: H12 becomes a SUM formula ---
$ws.Range("H12").Formula = "=SUM(H2:H11)"

# --- Format the new total cell: centered, custom font, white fill, medium right/bottom border ---
$total = $ws.Range("H12")
$total.Interior.Color = 16777215
$total.Borders.Item(10).Weight = -4138
$total.Borders.Item(9).Weight = -4138
$total.HorizontalAlignment = -4108
$total.VerticalAlignment = -4108
$total.Font.Name = "Century"
$total.Font.Family = 1

# --- Row 12 gets a slightly taller height to accommodate the thicker bottom border ---
$ws.Rows.Item(12).RowHeight = 16.5

# --- Update the active selection to match the saved workbook state ---
$ws.Range("H13").Select()
